$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O3").Value = 0.0005230903625488281
$ws.Range("O6").Value = 0.0009996891021728516
$ws.Range("O8").Value = 0.001005649566650391
$ws.Range("O11").Value = 0
$ws.Range("O13").Value = 0.001006603240966797
$ws.Range("O14").Value = 0.2026379108428955
$ws.Range("O15").Value = 0.001012563705444336
$ws.Range("O16").Value = 0.001085758209228516
$ws.Range("O17").Value = 0.001104116439819336
$ws.Range("O18").Value = 0.00201106071472168
$ws.Range("O19").Value = 0.001000404357910156
$ws.Range("O20").Value = 0.002087116241455078
$ws.Range("O21").Value = 0.0009996891021728516
$ws.Range("O23").Value = 0.001104593276977539
$ws.Range("O25").Value = 0.1101336479187012
$ws.Range("O26").Value = 0.1527211666107178
$ws.Range("O27").Value = 0.02155756950378418
$ws.Range("O28").Value = 0.001002073287963867
$ws.Range("O29").Value = 0.009074687957763672
$ws.Range("O30").Value = 0.004054546356201172
$ws.Range("O31").Value = 0.001998662948608398
$ws.Range("O33").Value = 0
$ws.Range("O34").Value = 12.7786328792572
$ws.Range("O35").Value = 0.02059102058410645
$ws.Range("O36").Value = 0.1852035522460938
$ws.Range("O37").Value = 0.3207833766937256
$ws.Range("O38").Value = 0.01584410667419434
$ws.Range("O39").Value = 0.01582050323486328
$ws.Range("O40").Value = 0.1353545188903809
$ws.Range("O42").Value = 0.1510148048400879
$ws.Range("O43").Value = 0.003079414367675781
$ws.Range("O44").Value = 0.1099669933319092
$ws.Range("O46").Value = 0.003000736236572266
$ws.Range("O47").Value = 0.001042366027832031
$ws.Range("O48").Value = 0.1147842407226562
$ws.Range("O51").Value = 0.007058143615722656
$ws.Range("O52").Value = 0
$ws.Range("O54").Value = 0.01011919975280762
$ws.Range("O55").Value = 0.06833338737487793
$ws.Range("O57").Value = 0.0111076831817627
$ws.Range("O61").Value = 0.001973152160644531
$ws.Range("O62").Value = 0.01300740242004395
$ws.Range("O63").Value = 0.0006785392761230469
$ws.Range("O64").Value = 0
$ws.Range("O65").Value = 0.09147787094116211
$ws.Range("O66").Value = 0.009242773056030273
$ws.Range("O68").Value = 0
$ws.Range("O69").Value = 0.01540732383728027
$ws.Range("O70").Value = 0.001875638961791992
$ws.Range("O71").Value = 0.02425336837768555
$ws.Range("O72").Value = 138.5756075382233
$ws.Range("O73").Value = 0.01060843467712402
$ws.Range("O75").Value = 0
$ws.Range("O76").Value = 0
$ws.Range("O77").Value = 0.01604676246643066
$ws.Range("O81").Value = 0.05000948905944824
$ws.Range("O85").Value = 0
$ws.Range("O86").Value = 0
$ws.Range("O87").Value = 0.007982969284057617
$ws.Range("O90").Value = 0.1330914497375488
$ws.Range("O96").Value = 0.08298754692077637
$ws.Range("O97").Value = 0
$ws.Range("O98").Value = 0.008511781692504883
$ws.Range("O100").Value = 0.07381963729858398
$ws.Range("O103").Value = 0.03437972068786621
$ws.Range("O104").Value = 3.656021595001221
$ws.Range("O105").Value = 0.02335095405578613
$ws.Range("O106").Value = 0.04840922355651855
$ws.Range("O107").Value = 0.03774404525756836
$ws.Range("O108").Value = 0.01321983337402344
$ws.Range("O109").Value = 0.01050257682800293
$ws.Range("O110").Value = 0.05618071556091309
$ws.Range("O111").Value = 0
$ws.Range("O112").Value = 0.02182388305664062
$ws.Range("O113").Value = 0.03330779075622559
$ws.Range("O114").Value = 0.01650452613830566
$ws.Range("O115").Value = 0.01251626014709473
$ws.Range("O116").Value = 0.0126183032989502
$ws.Range("O117").Value = 0.0210573673248291
$ws.Range("O118").Value = 0
$ws.Range("O120").Value = 0
$ws.Range("O121").Value = 0
$ws.Range("O122").Value = 0
$ws.Range("O126").Value = 0
$ws.Range("O129").Value = 0.001497268676757812
$ws.Range("O131").Value = 0.0003156661987304688
$ws.Range("O132").Value = 0
$ws.Range("O133").Value = 0
$ws.Range("O134").Value = 0
$ws.Range("O135").Value = 0.0005707740783691406
$ws.Range("O136").Value = 0.0009629726409912109
$ws.Range("O137").Value = 0
$ws.Range("O138").Value = 0.0009996891021728516
$ws.Range("O139").Value = 0
$ws.Range("O140").Value = 0
$ws.Range("O141").Value = 0
$ws.Range("O142").Value = 0.0161592960357666
$ws.Range("O143").Value = 0
$ws.Range("O144").Value = 0
$ws.Range("O145").Value = 0
$ws.Range("O146").Value = 0.001001119613647461
$ws.Range("O148").Value = 0.01062870025634766
$ws.Range("O149").Value = 0
$ws.Range("O150").Value = 0.2923462390899658
$ws.Range("O151").Value = 0.002306938171386719
$ws.Range("O153").Value = 0
$ws.Range("O159").Value = 0
$ws.Range("O160").Value = 0
$ws.Range("O161").Value = 0
$ws.Range("O164").Value = 0
$ws.Range("O165").Value = 0.001004219055175781
$ws.Range("O166").Value = 0.000919342041015625
$ws.Range("O170").Value = 0.001102209091186523
$ws.Range("O172").Value = 0.008905172348022461
$ws.Range("O174").Value = 0
$ws.Range("O175").Value = 0
$ws.Range("O176").Value = 0.005504846572875977
